$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 202, shifting existing rows 202:301 down to 203:302
$ws.Rows("202:202").Insert()

# Populate the newly inserted row 202 with the new weekly data record
$ws.Range("A202").Value = 6
$ws.Range("B202").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C202").Value = "Metropolitana"
$ws.Range("D202").Value = 44992
$ws.Range("E202").Value = 13
$ws.Range("F202").Value = 100112029
$ws.Range("G202").Value = "Orégano"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 44
$ws.Range("K202").Value = 16000
$ws.Range("L202").Value = 17000
$ws.Range("M202").Value = 16477
$ws.Range("N202").Value = "$/docena de atados"
$ws.Range("O202").Value = "Región Metropolitana"
$ws.Range("P202").Value = 5492
$ws.Range("Q202").Value = 3
$ws.Range("R202").Value = "Hortaliza"
